$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the longer task descriptions
$ws.Columns.Item(2).ColumnWidth = 78.92

# Row 39
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A39:C39").PasteSpecial(-4122) | Out-Null

# Row 40
$ws.Range("A40").Value = "Feb 03 10:00 to 11:00"
$ws.Range("B40").Value = "Created features: max event occur, max sublocation occur, min`nsublocation occur, alarm duration."
$ws.Range("C40").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A40:C40").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(40).RowHeight = 30

# Row 41
$ws.Range("A41").Value = "Feb 03 11:00 to 12:00"
$ws.Range("B41").Value = "Created features: day of week, day of month and month of year."
$ws.Range("C41").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A41:C41").PasteSpecial(-4122) | Out-Null

# Row 42
$ws.Range("A42").Value = "Feb 03 12:00 to 13:00"
$ws.Range("B42").Value = "Filled na values of max event occur with event close to mean of `nevents."
$ws.Range("C42").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A42:C42").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(42).RowHeight = 30

# Row 43
$ws.Range("A43").Value = "Feb 03 13:00 to 14:00"
$ws.Range("B43").Value = "Structurized problem statement. Understood input variables, output`nvariables and at what factors output should be evaluated. Understood`ncategorical attributes in dataset. Converting the categorical values to`nnumerical using label encoding method."
$ws.Range("C43").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A43:C43").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(43).RowHeight = 105

# Row 44
$ws.Range("A44").Value = "Feb 03 14:00 to 14:30"
$ws.Range("B44").Value = "Lunch"
$ws.Range("C44").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A44:C44").PasteSpecial(-4122) | Out-Null

# Row 45
$ws.Range("A45").Value = "Feb 03 14:00 to 15:00"
$ws.Range("B45").Value = "Applied label encoding to categorical columns."
$ws.Range("C45").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A45:C45").PasteSpecial(-4122) | Out-Null

# Row 46
$ws.Range("A46").Value = "Feb 03 15:00 to 16:00"
$ws.Range("B46").Value = "Applied multi output regression model. Generated train-test split.`nUsed train-test ratio as 70:30"
$ws.Range("C46").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A46:C46").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(46).RowHeight = 30

# Row 47
$ws.Range("A47").Value = "Feb 03 16:00 to 17:00"
$ws.Range("B47").Value = "Applied SVM algo. Checked model accuracy after it."
$ws.Range("C47").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A47:C47").PasteSpecial(-4122) | Out-Null

# Row 48
$ws.Range("A48").Value = "Feb 03 17:00 to 18:00"
$ws.Range("B48").Value = "Printed model accuracy of both the algorithms."
$ws.Range("C48").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A48:C48").PasteSpecial(-4122) | Out-Null

# Row 49
$ws.Range("A49").Value = "Feb 03 18:00 to 19:00"
$ws.Range("B49").Value = "Working on scalling because of results of predictions were not good.`nGeeting join conflicts while concatinating scaled and non scaled data.`nWorking on finding best hyperparametrs and also applying deep`nlearning techniques."
$ws.Range("C49").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A49:C49").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(49).RowHeight = 90

# Row 50
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A50:C50").PasteSpecial(-4122) | Out-Null

# Row 51
$ws.Range("A51").Value = "Feb 04 10:00 to 11:00"
$ws.Range("B51").Value = "Scaled input variable within various ranges every time. The results`nwere not great. Using KNN regresor for prediction."
$ws.Range("C51").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A51:C51").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(51).RowHeight = 30

# Row 52
$ws.Range("A52").Value = "Feb 04 11:00 to 12:00"
$ws.Range("B52").Value = "Aplied KNN regressor on dataset for prediction. Visualized RMSE from`nk equals 1 to 50."
$ws.Range("C52").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A52:C52").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(52).RowHeight = 45

# Row 53
$ws.Range("A53").Value = "Feb 04 12:00 to 13:00"
$ws.Range("B53").Value = "Found best hyperparameters for knn. Finding best hyperparameters`nfor random forest."
$ws.Range("C53").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A53:C53").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(53).RowHeight = 30

# Row 54
$ws.Range("A54").Value = "Feb 04 13:30 to 14:00"
$ws.Range("B54").Value = "Printed best hyperparameters of knn and random forest."
$ws.Range("C54").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A54:C54").PasteSpecial(-4122) | Out-Null
$ws.Range("A54").NumberFormat = "mmm-yy"
$ws.Range("A54").HorizontalAlignment = -4108
$ws.Range("A54").VerticalAlignment = -4108

# Row 55
$ws.Range("A55").Value = "Feb 04 14:00 to 15:00"
$ws.Range("B55").Value = "Generating sample predictions"
$ws.Range("C55").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A55:C55").PasteSpecial(-4122) | Out-Null

# Row 56
$ws.Range("A56").Value = "Feb 04 15:00 to 16:00"
$ws.Range("B56").Value = "Choosed random forest regressor for building model. Generated`nsample prediction on data."
$ws.Range("C56").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A56:C56").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(56).RowHeight = 30

# Row 57
$ws.Range("A57").Value = "Feb 04 16:00 to 17:00"
$ws.Range("B57").Value = "Undertanding tensorflow for modeling regresion problem. Solving an`nexample to understand prediction using tensorflow and keras."
$ws.Range("C57").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A57:C57").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(57).RowHeight = 45

# Row 58
$ws.Range("A58").Value = "Feb 04 17:00 to 18:00"
$ws.Range("B58").Value = "Understanding tensorflow for regression problems."
$ws.Range("C58").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A58:C58").PasteSpecial(-4122) | Out-Null

# Row 59
$ws.Range("A59").Value = "Feb 04 18:00 to 19:00"
$ws.Range("B59").Value = "Understand and applied neural to dataset."
$ws.Range("C59").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A59:C59").PasteSpecial(-4122) | Out-Null

# Row 60
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A60:C60").PasteSpecial(-4122) | Out-Null

# Row 61
$ws.Range("A61").Value = "Feb 05 10:00 to 11:00"
$ws.Range("B61").Value = "Modified some code of combined analysis. Added doc strings in`nfunctions."
$ws.Range("C61").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A61:C61").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(61).RowHeight = 30

# Row 62
$ws.Range("A62").Value = "Feb 05 11:00 to 12:00"
$ws.Range("B62").Value = "Modifying logic of imputing null values. Writen function which`nimputes null value with random value picked within range mean - std`nand mean + std."
$ws.Range("C62").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A62:C62").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(62).RowHeight = 60

# Row 63
$ws.Range("A63").Value = "Feb 05 12:00 to 13:00"
$ws.Range("B63").Value = "Tried every possible way to impute value using above function, but`nevery time got error."
$ws.Range("C63").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A63:C63").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(63).RowHeight = 30

# Row 64
$ws.Range("A64").Value = "Feb 05 13:00 to 13:30"
$ws.Range("B64").Value = "Done imputing na values using above function."
$ws.Range("C64").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A64:C64").PasteSpecial(-4122) | Out-Null

# Row 65
$ws.Range("A65").Value = "Feb 05 13:30 to 14:00"
$ws.Range("B65").Value = "Lunch"
$ws.Range("C65").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A65:C65").PasteSpecial(-4122) | Out-Null

# Row 66
$ws.Range("A66").Value = "Feb 05 14:00 to 15:00"
$ws.Range("B66").Value = "Corrected code for filling na using above funcion. Modefied code of`nfilling na of categorical columns."
$ws.Range("C66").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A66:C66").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(66).RowHeight = 30

# Row 67
$ws.Range("A67").Value = "Feb 05 15:00 to 16:00"
$ws.Range("B67").Value = "Understanding regresion for deep learning."
$ws.Range("C67").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A67:C67").PasteSpecial(-4122) | Out-Null

# Row 68
$ws.Range("A68").Value = "Feb 05 16:00 to 17:00"
$ws.Range("B68").Value = "Installed tensorflow and doing tensorflow regression example"
$ws.Range("C68").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A68:C68").PasteSpecial(-4122) | Out-Null

# Row 69
$ws.Range("A69").Value = "Feb 05 17:00 to 18:00"
$ws.Range("B69").Value = "Tensorflow caused problem while importing or loading"
$ws.Range("C69").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A69:C69").PasteSpecial(-4122) | Out-Null

# Row 70
$ws.Range("A70").Value = "Feb 05 18:00 to 19:00"
$ws.Range("B70").Value = "Solved sample examples of temsorflow, caused error while programming. Class"
$ws.Range("C70").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A70:C70").PasteSpecial(-4122) | Out-Null

# Row 71
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A71:C71").PasteSpecial(-4122) | Out-Null

# Row 72
$ws.Range("A72").Value = "Feb 06 10:00 to 11:00"
$ws.Range("B72").Value = "Checking colinerity and covariance of variables."
$ws.Range("C72").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A72:C72").PasteSpecial(-4122) | Out-Null

# Row 73
$ws.Range("A73").Value = "Feb 06 11:00 to 12:00"
$ws.Range("B73").Value = "Plotter scater plot of some variables and checked colinearity. Found the reason why`ntotal output count results were poor because of no proper relationship with other`nvariables."
$ws.Range("C73").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A73:C73").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(73).RowHeight = 60

# Row 74
$ws.Range("A74").Value = "Feb 06 12:00 to 13:00"
$ws.Range("B74").Value = "Documented observations and solutions."
$ws.Range("C74").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A74:C74").PasteSpecial(-4122) | Out-Null

# Row 75
$ws.Range("A75").Value = "Feb 06 13:00 to 13:30"
$ws.Range("B75").Value = "Practicing simple neural network"
$ws.Range("C75").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A75:C75").PasteSpecial(-4122) | Out-Null

# Row 76
$ws.Range("A76").Value = "Feb 06 13:30 to 14:00"
$ws.Range("B76").Value = "Lunch"
$ws.Range("C76").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A76:C76").PasteSpecial(-4122) | Out-Null

# Row 77
$ws.Range("A77").Value = "Feb 06 14:00 to 15:00"
$ws.Range("B77").Value = "Practicing neural network, implemented neural network class with fit, activation and`ntraining functionalities."
$ws.Range("C77").Value = "Infimetrics"
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A77:C77").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(77).RowHeight = 45

# Row 78
$ws.Range("A78").Value = "Feb 06 15:00 to 16:00"
$ws.Range("B78").Value = "Implemented simple neural network program. Added functionality for tanh activation."
$ws.Range("C78").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A78:C78").PasteSpecial(-4122) | Out-Null

# Row 79
$ws.Range("A79").Value = "Feb 06 16:00 to 17:00"
$ws.Range("B79").Value = "Implemented rmse and mse accuracy functionality."
$ws.Range("C79").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A79:C79").PasteSpecial(-4122) | Out-Null

# Row 80
$ws.Range("A80").Value = "Feb 06 17:00 to 18:00"
$ws.Range("B80").Value = "Implemented multi layer perceptron. Added functionalities of dot product and activation"
$ws.Range("C80").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A80:C80").PasteSpecial(-4122) | Out-Null

# Row 81
$ws.Range("A81").Value = "Feb 06 18:00 to 19:00"
$ws.Range("B81").Value = "Implemented back propagation for MLP."
$ws.Range("C81").Value = "Infimetrics"
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A81:C81").PasteSpecial(-4122) | Out-Null

# Update the view to match where the author ended up working
$ws.Range("B77").Select() | Out-Null
